# Adds LOC_PROD, LOC_TEST, LOC_TOT, Q1, Q2, Q3 columns (P:U) to the Tasks1&2 sheet,
# mirroring the header/data-row styling of the existing I:O "code metric" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header row (P1:U1), with the same style as the existing O1 header ----
$headers = @("LOC_PROD", "LOC_TEST", "LOC_TOT", "Q1", "Q2", "Q3")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 16 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Copy the header style (bold font + bottom border + numFmt) from O1 onto P1:U1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:U1").PasteSpecial(-4122) | Out-Null

# ---- New data values (P2:U19), matching the style used by the other numeric columns ----
$values = @{
    2  = @(102, 60, 162, 5, 3, 3)
    3  = @(146, 81, 207, 3, 3, 3)
    4  = @(90,  64, 154, 4, 2, 3)
    5  = @(121, 57, 178, 4, 3, 2)
    6  = @(93, 102, 195, 5, 4, 4)
    7  = @(164, 96, 260, 5, 3, 3)
    8  = @(90,  67, 157, 3, 3, 3)
    9  = @(111, 74, 185, 4, 3, 4)
    10 = @(82,   8, 100, 4, 3, 3)
    11 = @(86,  27, 113, 2, 2, 1)
    12 = @(90,  45, 135, 5, 4, 4)
    13 = @(74,   6,  80, 2, 1, 1)
    14 = @(85,  72, 157, 5, 5, 5)
    15 = @(132, 71, 203, 4, 5, 5)
    16 = @(84,   8,  92, 5, 3, 3)
    17 = @(87,  45, 129, 5, 3, 2)
    18 = @(66,   8,  74, 2, 2, 3)
    19 = @(99,  16, 115, 2, 1, 2)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = 16 + $i
        $ws.Cells.Item($row, $col).Value = $rowVals[$i]
    }
}

# Copy the plain numeric-column style (used by F2:H19 etc.) onto the new data block
$ws.Range("F2").Copy() | Out-Null
$ws.Range("P2:U19").PasteSpecial(-4122) | Out-Null

# ---- Column widths: extend the existing 10.7109375-wide formatting to P:U ----
$ws.Range("P1:U1").EntireColumn.ColumnWidth = $ws.Range("A1").ColumnWidth
$excel.CutCopyMode = $false

# ---- Selection, as left by the author after the edit ----
$ws.Range("Q22").Select() | Out-Null
